# Auto-generated edit script applying the cryptos.xlsx data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '71.494.61'
$ws.Range("E2").Value = '  -2.19%  '
$ws.Range("D3").Value = '3.886.53'
$ws.Range("E3").Value = '  -2.57%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '607.12'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.59%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '172.88'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +5.38%  '
$ws.Range("E7").Value = '  -1.97%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.754'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.21%  '
$ws.Range("E10").Value = '  +5.73%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '54.17'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.71%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000325'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.26%  '
$ws.Range("E13").Value = '  +4.73%  '
$ws.Range("D14").Value = '4.501.44'
$ws.Range("E14").Value = '  -2.75%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '21.26'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.39%  '
$ws.Range("D16").Value = '3.879.65'
$ws.Range("E16").Value = '  -2.89%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.01'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.03%  '
$ws.Range("E18").Value = '  -4.23%  '
$ws.Range("E19").Value = '  -2.17%  '
$ws.Range("D20").Value = '71.263.02'
$ws.Range("E20").Value = '  -2.10%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '442.35'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.29%  '
$ws.Range("E22").Value = '  -0.78%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '95.07'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.98%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.32'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.13%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '13.97'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.23%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.89'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.23%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '4.05'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -6.75%  '
$ws.Range("E28").Value = '  +0.07%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.56'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.96%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.79'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +11.85%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '35.35'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.18%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '13.61'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.36%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '48.16'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.76%  '
$ws.Range("E34").Value = '  -3.76%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0000100'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +10.69%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '69.73'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.02%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '637.97'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.30%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.444'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.07%  '
$ws.Range("E39").Value = '  +0.83%  '
$ws.Range("E40").Value = '  +0.06%  '
$ws.Range("E41").Value = '  -0.08%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.29'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.49%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.91'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +9.61%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.20'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +19.69%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0474'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.44%  '
$ws.Range("E46").Value = '  -4.42%  '
$ws.Range("B47").Value = 'Stellar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.145'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.05%  '
$ws.Range("B48").Value = 'WEMIXToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.91'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -12.61%  '
$ws.Range("D49").Value = '2.926.66'
$ws.Range("E49").Value = '  +0.59%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.27'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.97%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.000277'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.64%  '
